$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 890, pushing all existing rows (old 890-931)
# down to become rows 892-933.
$ws.Rows.Item(890).Insert()
$ws.Rows.Item(890).Insert()

# New row 890: 2026/02/26, 木, 20, 201
$ws.Range("A890").NumberFormat = "@"
$ws.Range("A890").Value = "2026/02/26"
$ws.Range("A890").Style = "Normal"
$ws.Range("B890").NumberFormat = "@"
$ws.Range("B890").Value = "木"
$ws.Range("B890").Style = "Normal"
$ws.Range("C890").Value = 20
$ws.Range("D890").Value = 201

# New row 891: 2026/02/26, 木, 22, 201
$ws.Range("A891").NumberFormat = "@"
$ws.Range("A891").Value = "2026/02/26"
$ws.Range("A891").Style = "Normal"
$ws.Range("B891").NumberFormat = "@"
$ws.Range("B891").Value = "木"
$ws.Range("B891").Style = "Normal"
$ws.Range("C891").Value = 22
$ws.Range("D891").Value = 201

$ws.Range("A1").Select()
